$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Original layout:
#   A10 = "Soll Kriterien:" (bold header style)
#   A11 = "Loeschen von Interessenprofilen" (regular style)
# New layout:
#   A9  = "Loeschen von Interessenprofilen" (regular style)
#   A10 = empty
#   A11 = "Soll Kriterien:" (bold header style)

# Move A11 -> A9 (Cut preserves both value and formatting)
$ws.Range("A11").Cut($ws.Range("A9"))

# Move A10 -> A11 (Cut preserves both value and formatting)
$ws.Range("A10").Cut($ws.Range("A11"))

# The cut-from cell keeps its old formatting; reset it back to the
# workbook's default (unformatted) style so A10 ends up fully empty.
$ws.Range("A10").Style = "Standard"

$ws.Range("A11").Select()
